# DOMA-3100 add formatter convert to number for some colomns
#
# For the ticket analytics export template, append the ":formatN()" formatter
# to the placeholder text of the numeric columns (processing, completed,
# canceled, deferred, closed, new_or_reopened) in the two data rows, and
# switch those cells' number format from text ("@") to a plain number ("0")
# so the exported values are treated/rendered as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C..H hold: processing, completed, canceled, deferred, closed, new_or_reopened
$cols = @("C", "D", "E", "F", "G", "H")

foreach ($row in 2, 3) {
    foreach ($col in $cols) {
        $cell = $ws.Range($col + $row)
        $text = $cell.Value2
        # Insert the ":formatN()" formatter right before the closing brace,
        # e.g. "{d.tickets[i].processing}" -> "{d.tickets[i].processing:formatN()}"
        $cell.Value2 = $text.Substring(0, $text.Length - 1) + ":formatN()}"
        $cell.NumberFormat = "0"
    }
}
